$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Clear project-name / PAT values in row 2, keep the hyperlink cell text (A2)
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = $null

# Row 3 becomes completely empty, including removing the hyperlink on A3
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("A3").Value = $null
$ws.Range("B3").Value = $null
$ws.Range("C3").Value = $null

# Move the active selection to C2
$ws.Range("C2").Select()
